# Apply cryptos list update (price/volume refresh + Polkadot/WrappedEther row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.495.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.670.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.96%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.02"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5148"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.74%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06470"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2577"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.69%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07666"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.56%  "

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.344"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.73%  "

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.669.58"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.901.17"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5576"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8066"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.80"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.529.68"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.009"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.81"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.438"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.13"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.908"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.65%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.91"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.737"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1166"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.013"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.77"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.96%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05224"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.262"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.379"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.221"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.582"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.770"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.51%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9259"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5746"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.164.10"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +11.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01600"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8495"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.43%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.643"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.36"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.810.51"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈111"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4493"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.08"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.99%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.975"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05131"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.11%  "
